# The DataFrame was originally written to Excel without passing
# index=False, so an extra "index" column ended up in column A.
# This fix removes that leftover index column, shifting the real
# data (First Name, Email, Subscription Date) one column to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire first column (the stray pandas index values),
# which shifts columns B:D left to A:C.
$ws.Columns.Item(1).Delete()
